# Updates cached market-board price/profit figures across several Leve sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), mirroring a scheduled data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5802.222
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 6152.5
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 6152.5
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -7400.5
$ws.Range("H65").Value = 5802.222
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 6152.5
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 30762.5
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -37002.5
$ws.Range("H98").Value = 55556572
$ws.Range("I98").Value = 55556572
$ws.Range("K98").Value = 55556572
$ws.Range("M98").Value = -55555074
$ws.Range("H122").Value = 55556572
$ws.Range("I122").Value = 55556572
$ws.Range("K122").Value = 166669716
$ws.Range("M122").Value = -166667266
$ws.Range("H133").Value = 69999.5
$ws.Range("J133").Value = 69999.5
$ws.Range("L133").Value = 69999.5
$ws.Range("N133").Value = -80119.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 627.12
$ws.Range("I2").Value = 545.6842
$ws.Range("K2").Value = 545.6842
$ws.Range("M2").Value = -432.6842
$ws.Range("H32").Value = 8337559
$ws.Range("I32").Value = 10205358
$ws.Range("K32").Value = 10205358
$ws.Range("M32").Value = -10205071
$ws.Range("H61").Value = 20277122
$ws.Range("I61").Value = 14712427
$ws.Range("J61").Value = 83343670
$ws.Range("K61").Value = 14712427
$ws.Range("L61").Value = 83343670
$ws.Range("M61").Value = -14712215
$ws.Range("N61").Value = -83344094
$ws.Range("H74").Value = 13005286
$ws.Range("I74").Value = 16669137
$ws.Range("J74").Value = 2013735.2
$ws.Range("K74").Value = 16669137
$ws.Range("L74").Value = 2013735.2
$ws.Range("M74").Value = -16668263
$ws.Range("N74").Value = -2015483.2
$ws.Range("H77").Value = 13005286
$ws.Range("I77").Value = 16669137
$ws.Range("J77").Value = 2013735.2
$ws.Range("K77").Value = 83345685
$ws.Range("L77").Value = 10068676
$ws.Range("M77").Value = -83341317
$ws.Range("N77").Value = -10077412
$ws.Range("H116").Value = 627.12
$ws.Range("I116").Value = 545.6842
$ws.Range("K116").Value = 545.6842
$ws.Range("M116").Value = 1748.3158
$ws.Range("H132").Value = 13111
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 17916.5
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 53749.5
$ws.Range("M132").Value = -7970
$ws.Range("N132").Value = -58809.5
$ws.Range("H136").Value = 20277122
$ws.Range("I136").Value = 14712427
$ws.Range("J136").Value = 83343670
$ws.Range("K136").Value = 44137281
$ws.Range("L136").Value = 250031010
$ws.Range("M136").Value = -44134731
$ws.Range("N136").Value = -250036110

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 627.12
$ws.Range("I3").Value = 545.6842
$ws.Range("K3").Value = 545.6842
$ws.Range("M3").Value = -431.6842
$ws.Range("H95").Value = 38675
$ws.Range("J95").Value = 38675
$ws.Range("L95").Value = 38675
$ws.Range("N95").Value = -44167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7252.091
$ws.Range("J86").Value = 6924.4287
$ws.Range("L86").Value = 6924.4287
$ws.Range("N86").Value = -9170.4287
$ws.Range("H89").Value = 7252.091
$ws.Range("J89").Value = 6924.4287
$ws.Range("L89").Value = 34622.14350000001
$ws.Range("N89").Value = -45854.14350000001
$ws.Range("H99").Value = 3825.2222
$ws.Range("I99").Value = 3546
$ws.Range("J99").Value = 4048.6
$ws.Range("K99").Value = 3546
$ws.Range("L99").Value = 4048.6
$ws.Range("M99").Value = -2048
$ws.Range("N99").Value = -7044.6
$ws.Range("H125").Value = 980000
$ws.Range("J125").Value = 980000
$ws.Range("L125").Value = 980000
$ws.Range("N125").Value = -984920
$ws.Range("H126").Value = 3825.2222
$ws.Range("I126").Value = 3546
$ws.Range("J126").Value = 4048.6
$ws.Range("K126").Value = 10638
$ws.Range("L126").Value = 12145.8
$ws.Range("M126").Value = -8168
$ws.Range("N126").Value = -17085.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 84474.5
$ws.Range("J37").Value = 84474.5
$ws.Range("L37").Value = 253423.5
$ws.Range("N37").Value = -253647.5
$ws.Range("H117").Value = 1721.2222
$ws.Range("J117").Value = 1721.2222
$ws.Range("L117").Value = 5163.6666
$ws.Range("N117").Value = -12047.6666
$ws.Range("H131").Value = 7593.6123
$ws.Range("J131").Value = 7745.0415
$ws.Range("L131").Value = 23235.1245
$ws.Range("N131").Value = -33315.12450000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6326.222
$ws.Range("I102").Value = 11975
$ws.Range("J102").Value = 4712.2856
$ws.Range("K102").Value = 11975
$ws.Range("L102").Value = 4712.2856
$ws.Range("M102").Value = -10353
$ws.Range("N102").Value = -7956.2856
$ws.Range("H132").Value = 30310540
$ws.Range("I132").Value = 55559420
$ws.Range("J132").Value = 11884.2
$ws.Range("K132").Value = 166678260
$ws.Range("L132").Value = 35652.60000000001
$ws.Range("M132").Value = -166675730
$ws.Range("N132").Value = -40712.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 73446
$ws.Range("I7").Value = 3139.4
$ws.Range("K7").Value = 3139.4
$ws.Range("M7").Value = -3027.4
$ws.Range("H40").Value = 4718.5454
$ws.Range("I40").Value = 3978.6667
$ws.Range("J40").Value = 5230.769
$ws.Range("K40").Value = 3978.6667
$ws.Range("L40").Value = 5230.769
$ws.Range("M40").Value = -3842.6667
$ws.Range("N40").Value = -5502.769
$ws.Range("H61").Value = 2021.6
$ws.Range("I61").Value = 1876
$ws.Range("J61").Value = 2240
$ws.Range("K61").Value = 1876
$ws.Range("L61").Value = 2240
$ws.Range("M61").Value = -1674
$ws.Range("N61").Value = -2644
$ws.Range("H68").Value = 2850.1667
$ws.Range("I68").Value = 2800.2
$ws.Range("J68").Value = 3100
$ws.Range("K68").Value = 2800.2
$ws.Range("L68").Value = 3100
$ws.Range("M68").Value = -2051.2
$ws.Range("N68").Value = -4598
$ws.Range("H71").Value = 2850.1667
$ws.Range("I71").Value = 2800.2
$ws.Range("J71").Value = 3100
$ws.Range("K71").Value = 14001
$ws.Range("L71").Value = 15500
$ws.Range("M71").Value = -10257
$ws.Range("N71").Value = -22988
$ws.Range("H113").Value = 2021.6
$ws.Range("I113").Value = 1876
$ws.Range("J113").Value = 2240
$ws.Range("K113").Value = 1876
$ws.Range("L113").Value = 2240
$ws.Range("M113").Value = 294
$ws.Range("N113").Value = -6580
$ws.Range("H122").Value = 5451.8823
$ws.Range("I122").Value = 4998.96
$ws.Range("K122").Value = 14996.88
$ws.Range("M122").Value = -12546.88
$ws.Range("H126").Value = 73446
$ws.Range("I126").Value = 3139.4
$ws.Range("K126").Value = 9418.200000000001
$ws.Range("M126").Value = -6948.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 811.15625
$ws.Range("I107").Value = 852.4583
$ws.Range("K107").Value = 2557.3749
$ws.Range("M107").Value = -637.3748999999998
$ws.Range("H113").Value = 688.7059
$ws.Range("I113").Value = 507.3
$ws.Range("K113").Value = 1521.9
$ws.Range("M113").Value = 648.0999999999999
$ws.Range("H126").Value = 1308.0476
$ws.Range("I126").Value = 1188.45
$ws.Range("J126").Value = 3700
$ws.Range("K126").Value = 3565.35
$ws.Range("L126").Value = 11100
$ws.Range("M126").Value = -1095.35
$ws.Range("N126").Value = -16040
$ws.Range("H135").Value = 156485.8
$ws.Range("I135").Value = 45000
$ws.Range("K135").Value = 45000
$ws.Range("M135").Value = -39930

